$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 102; this shifts the former rows
# 102-181 down to 103-182 (and the sheet dimension grows to A1:R182
# automatically).
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new record.
$ws.Range("A102").Value2 = 2
$ws.Range("B102").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C102").Value2 = "Coquimbo"
$ws.Range("D102").Value2 = 44741
$ws.Range("E102").Value2 = 4
$ws.Range("F102").Value2 = 100112043
$ws.Range("G102").Value2 = "Pepino ensalada"
$ws.Range("H102").Value2 = "Sin especificar"
$ws.Range("I102").Value2 = "Primera"
$ws.Range("J102").Value2 = 500
$ws.Range("K102").Value2 = 15000
$ws.Range("L102").Value2 = 17000
$ws.Range("M102").Value2 = 16000
$ws.Range("N102").Value2 = "$/caja 60 unidades"
$ws.Range("O102").Value2 = "Provincia de Limarí"
$ws.Range("P102").Value2 = 267
$ws.Range("Q102").Value2 = 60
$ws.Range("R102").Value2 = "Hortaliza"
